$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = theta_se values (previously all "(nan)")
$ws.Range("B4").Value = "(0.13)"
$ws.Range("C4").Value = "(0.2)"
$ws.Range("D4").Value = "(0.66)"
$ws.Range("E4").Value = "(2.1)"
$ws.Range("F4").Value = "(1.79)"
$ws.Range("G4").Value = "(0.61)"
$ws.Range("H4").Value = "(3.09)"
$ws.Range("I4").Value = "(0.33)"
$ws.Range("J4").Value = "(3.38)"
$ws.Range("K4").Value = "(0.47)"
$ws.Range("L4").Value = "(1.49)"

# Row 6 = lambda_se values (previously all "(nan)")
$ws.Range("B6").Value = "(0.45)"
$ws.Range("C6").Value = "(0.22)"
$ws.Range("D6").Value = "(0.64)"
$ws.Range("E6").Value = "(1.17)"
$ws.Range("F6").Value = "(1.42)"
$ws.Range("G6").Value = "(0.61)"
$ws.Range("H6").Value = "(2.43)"
$ws.Range("I6").Value = "(0.02)"
$ws.Range("J6").Value = "(1.37)"
$ws.Range("K6").Value = "(1.11)"
$ws.Range("L6").Value = "(1.38)"
